$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra optimizer columns (G:J) entirely - the study was narrowed
# down to just gdpa-pd / gdpa-adam / gdpa-noise / gdpa-direct.
$ws.Range("G:J").Delete() | Out-Null

# Rename the remaining headers to the shorter optimizer names.
$ws.Cells.Item(1, 3).Value2 = "gdpa-pd"
$ws.Cells.Item(1, 4).Value2 = "gdpa-adam"
$ws.Cells.Item(1, 5).Value2 = "gdpa-noise"
$ws.Cells.Item(1, 6).Value2 = "gdpa-direct"

# Refresh the results matrix (C:F) with the new run's numbers. Only cells
# whose value actually changed vs. the prior run are touched.
$ws.Cells.Item(7, 3).Value2 = 1.08
$ws.Cells.Item(7, 6).Value2 = 1.1
$ws.Cells.Item(8, 3).Value2 = 1.18
$ws.Cells.Item(8, 4).Value2 = 1.14
$ws.Cells.Item(8, 5).Value2 = 1.16
$ws.Cells.Item(9, 3).Value2 = 1.2
$ws.Cells.Item(9, 4).Value2 = 1.22
$ws.Cells.Item(9, 5).Value2 = 1.2
$ws.Cells.Item(9, 6).Value2 = 1.2
$ws.Cells.Item(10, 3).Value2 = 1.44
$ws.Cells.Item(10, 5).Value2 = 1.42
$ws.Cells.Item(10, 6).Value2 = 1.42
$ws.Cells.Item(11, 3).Value2 = 1.86
$ws.Cells.Item(11, 4).Value2 = 1.68
$ws.Cells.Item(11, 6).Value2 = 1.76
$ws.Cells.Item(12, 3).Value2 = 2.56
$ws.Cells.Item(12, 4).Value2 = 2.02
$ws.Cells.Item(12, 5).Value2 = 3.36
$ws.Cells.Item(12, 6).Value2 = 2.88
$ws.Cells.Item(13, 3).Value2 = 3.5
$ws.Cells.Item(13, 4).Value2 = 3.44
$ws.Cells.Item(13, 5).Value2 = 3.54
$ws.Cells.Item(13, 6).Value2 = 4.44
$ws.Cells.Item(14, 3).Value2 = 7
$ws.Cells.Item(14, 4).Value2 = 4.9375
$ws.Cells.Item(14, 5).Value2 = 4.319148936170213
$ws.Cells.Item(14, 6).Value2 = 3.765957446808511
$ws.Cells.Item(15, 3).Value2 = 4.555555555555555
$ws.Cells.Item(15, 4).Value2 = 3.608695652173913
$ws.Cells.Item(15, 5).Value2 = 9.088888888888889
$ws.Cells.Item(15, 6).Value2 = 6.642857142857143
$ws.Cells.Item(16, 3).Value2 = 4.463414634146342
$ws.Cells.Item(16, 4).Value2 = 5.8
$ws.Cells.Item(16, 5).Value2 = 8.512195121951219
$ws.Cells.Item(16, 6).Value2 = 11.1025641025641
$ws.Cells.Item(17, 3).Value2 = 11.2051282051282
$ws.Cells.Item(17, 4).Value2 = 10.41025641025641
$ws.Cells.Item(17, 5).Value2 = 13
$ws.Cells.Item(17, 6).Value2 = 13.46666666666667
$ws.Cells.Item(18, 3).Value2 = 8.857142857142858
$ws.Cells.Item(18, 4).Value2 = 14.06896551724138
$ws.Cells.Item(18, 5).Value2 = 22.125
$ws.Cells.Item(18, 6).Value2 = 16.65217391304348
$ws.Cells.Item(19, 3).Value2 = 16.36842105263158
$ws.Cells.Item(19, 4).Value2 = 20.55
$ws.Cells.Item(19, 5).Value2 = 20.85714285714286
$ws.Cells.Item(19, 6).Value2 = 18.91666666666667
$ws.Cells.Item(20, 3).Value2 = 11.6
$ws.Cells.Item(20, 4).Value2 = 13.8
$ws.Cells.Item(20, 5).Value2 = 18.71428571428572
$ws.Cells.Item(20, 6).Value2 = 46.28571428571428
$ws.Cells.Item(21, 3).Value2 = 8
$ws.Cells.Item(21, 4).Value2 = 29.25
$ws.Cells.Item(21, 5).Value2 = 28
$ws.Cells.Item(21, 6).Value2 = 47
